$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.298.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.411.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.840.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.188.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.415.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "326.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.173"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "325.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "149.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.44%  "
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
